# Row 3 is an exact duplicate of row 2's data (the same scraped match row
# appears twice in the source JSON/CSV this sheet was generated from).
# Copy the whole row 2 range into row 3 so every cell keeps its original
# text type and exact characters (e.g. the trailing non-breaking space in
# the batsman name, and numeric-looking values like "53.84" staying text)
# instead of being re-typed and risking reinterpretation as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:K2").Copy($ws.Range("A3:K3"))
